# Reformulación de Prez, Moodle y clase de test para identificar a su
# recurso por su ID: la hoja "Recursos" gana una columna "IDRecurso" con
# el identificador numérico del recurso de Moodle, y los alias/contextos
# se reescriben para reflejar los nuevos recursos de ejemplo.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Recursos")

# --- Insert a new column B ("IDRecurso") before the existing "Alias" column ---
$ws.Range("B:B").Insert()

# --- Header row ---
$ws.Range("A1").Value = "Contexto del evento"
$ws.Range("B1").Value = "IDRecurso"
$ws.Range("C1").Value = "Alias"
$ws.Range("D1").Value = "Excluido"

# --- Data rows: Contexto del evento | IDRecurso | Alias ---
$ws.Cells.Item(2, 1).Value = "Foro: Noticias de clase"
$ws.Cells.Item(2, 2).Value = 5000
$ws.Cells.Item(2, 3).Value = "Foro: Noticias de clase"

$ws.Cells.Item(3, 1).Value = "Carpeta: Exámenes"
$ws.Cells.Item(3, 2).Value = 5002
$ws.Cells.Item(3, 3).Value = "Carpeta: Exámenes"

$ws.Cells.Item(4, 1).Value = "Carpeta: Recursos del Alumnado"
$ws.Cells.Item(4, 2).Value = 5011
$ws.Cells.Item(4, 3).Value = "Carpeta: Recursos del Alumnado"

$ws.Cells.Item(5, 1).Value = "Carpeta: Recursos del Alumnado"
$ws.Cells.Item(5, 2).Value = 5012
$ws.Cells.Item(5, 3).Value = "Carpeta: Recursos del Alumnado"

$ws.Cells.Item(6, 1).Value = "Carpeta: Papeleo"
$ws.Cells.Item(6, 2).Value = 5013
$ws.Cells.Item(6, 3).Value = "Carpeta: Papeleo"

$ws.Cells.Item(7, 1).Value = "Carpeta: Recursos del Alumnado"
$ws.Cells.Item(7, 2).Value = 5014
$ws.Cells.Item(7, 3).Value = "Carpeta: Recursos del Alumnado"

$ws.Cells.Item(8, 1).Value = "Tarea: Entrega inicial"
$ws.Cells.Item(8, 2).Value = 5015
$ws.Cells.Item(8, 3).Value = "Tarea: Entrega inicial"

$ws.Cells.Item(9, 1).Value = "Carpeta: Entrega inicial"
$ws.Cells.Item(9, 2).Value = 5016
$ws.Cells.Item(9, 3).Value = "Carpeta: Entrega inicial"

# The old "Excluido" value ("NO") in row 2 is no longer applicable to the
# new first data row, so make sure column D stays empty for the data rows.
$ws.Range("D2:D9").ClearContents()

# --- Column widths / layout for the reshaped sheet ---
$ws.Range("A:A").ColumnWidth = 30.7109375
$ws.Range("B:B").ColumnWidth = 0
$ws.Range("B:B").EntireColumn.Hidden = $true
$ws.Range("C:C").ColumnWidth = 30.7109375

Write-Host "Recursos sheet reshaped"
